$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.998.12"
$ws.Range("E2").Value = "  +6.75%  "

$ws.Range("D3").Value = "3.353.92"
$ws.Range("E3").Value = "  +3.12%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'413.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.71%  "

$ws.Range("D6").Value = "'112.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "

$ws.Range("E7").Value = "  +3.88%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.641"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.43%  "

$ws.Range("D10").Value = "'39.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").Value = "'0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.35%  "

$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "3.889.65"
$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").Value = "'8.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.36%  "

$ws.Range("D15").Value = "'19.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.07%  "

$ws.Range("D16").Value = "3.325.67"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("D18").Value = "60.865.12"
$ws.Range("E18").Value = "  +6.74%  "

$ws.Range("D19").Value = "'10.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.25%  "

$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("E21").Value = "  +5.10%  "

$ws.Range("D22").Value = "'13.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'303.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("D24").Value = "'75.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.72%  "

$ws.Range("E25").Value = "  +1.87%  "

$ws.Range("D26").Value = "'28.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.85%  "

$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.52%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.52%  "

$ws.Range("E31").Value = "  +23.54%  "

$ws.Range("E32").Value = "  +4.79%  "

$ws.Range("D33").Value = "'11.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.07%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "'39.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.89%  "

$ws.Range("D36").Value = "'0.0509"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.93%  "

$ws.Range("D37").Value = "'52.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "

$ws.Range("D38").Value = "'3.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.91%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").Value = "'3.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("D41").Value = "'0.300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.82%  "

$ws.Range("D42").Value = "'137.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("E43").Value = "  +2.77%  "

$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("D45").Value = "'3.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").Value = "'16.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("E47").Value = "  +8.70%  "

$ws.Range("D48").Value = "'22.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("D49").Value = "2.187.92"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").Value = "'2.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("E51").Value = "  -2.22%  "
